# HungerRewards.xlsx edit: remove negative reward values, replace with "n"
# and delete the empty second sheet ("Sheet1").

$wb = $excel.ActiveWorkbook

# --- 1. Delete the empty extra sheet named "Sheet1" ---
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Sheet1") {
        $sheet.Delete()
    }
}

# --- 2. Update the data table on Tabelle1 ---
$ws = $wb.Worksheets.Item("Tabelle1")

# Rows 2-6: both C and D become "n"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = "n"
    $ws.Cells.Item($r, 4).Value = "n"
}

# Rows 7-16: only D becomes "n"
for ($r = 7; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = "n"
}

# --- 3. Update selection to D24 ---
$ws.Range("D24").Select()
